$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column C values for rows 13 to 43: change "TATA 1" -> "TATA 2"
for ($r = 13; $r -le 43; $r++) {
    $ws.Cells.Item($r, 3).Value = "TATA 2"
}

# Update the view: scroll/selection state
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 24
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C11:C43").Select()
